$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Alice, 25, New York, 50000
$ws.Range("E4").Value = 25
$ws.Range("F4").Value = "New York"

# Row 5: Bob, 30, Los Angeles, 60000
$ws.Range("E5").Value = 30
$ws.Range("F5").Value = "Los Angeles"

# Row 6: Charlie, (age empty), Chicago, 70000
$ws.Range("D6").Value = "Charlie"
$ws.Range("E6").ClearContents()
$ws.Range("F6").Value = "Chicago"
$ws.Range("G6").Value = 70000

# Row 7: David, 35, Houston, 80000
$ws.Range("D7").Value = "David"
$ws.Range("E7").Value = 35
$ws.Range("F7").Value = "Houston"
$ws.Range("G7").Value = 80000

# Row 8: Eve, 40, (city empty), 90000
$ws.Range("D8").Value = "Eve"
$ws.Range("E8").Value = 40
$ws.Range("F8").ClearContents()
$ws.Range("G8").Value = 90000

# Row 9: clear entirely
$ws.Range("C9:G9").ClearContents()

# Selection
$ws.Range("J13").Select()
